# Generate Report for Archive
#
# The localization status for the in-flight files moved on from
# "Ready for handoff" to "In Translation" - update every sheet that
# surfaces that status column (the Overview roll-up plus each
# per-locale report) and re-fit the now-shorter status columns the
# way Excel does after a content refresh.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New width (in characters) for the status column once it holds the
# shorter "In Translation" text instead of "Ready for handoff".
$statusColWidth = 12.5

# --- Overview sheet: status is mirrored per-locale in columns E (zh-cn)
#     and F (de-de), rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- Per-locale detail sheets: status lives in column C, rows 2-4 ---
$localeSheets = @("zh-cn", "de-de")
foreach ($sheetName in $localeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2:C4").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = $statusColWidth
}
